$d = $word.ActiveDocument

# --- Pass 1: walk the document's paragraphs in order and figure out which
#     "List Number" paragraphs still need a numbered-list (w:numPr) applied.
#     Word already numbers the "User chooses menu..." / "The system displays
#     a tab..." steps in each flow (ListType 3); the step(s) that follow
#     them in the same flow share that same list (numId) but are missing
#     w:numPr. We record, for every un-numbered candidate paragraph, the
#     numId of the nearest preceding numbered paragraph.
$paras = $d.Paragraphs
$count = $paras.Count

$targets = New-Object System.Collections.ArrayList
$lastNumId = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $lt = $p.Range.ListFormat.ListType
    if ($lt -eq 3) {
        $lastNumId = $p.Range.ListFormat.List.ListID
    } elseif ($lt -eq 0) {
        $styleName = $p.Style.NameLocal
        if ($styleName -eq "List Number" -and $lastNumId -ne -1) {
            $li = $p.Format.LeftIndent
            $is1296 = ($li -gt 64.7 -and $li -lt 64.9)
            $is1440 = ($li -gt 71.9 -and $li -lt 72.1)
            if ($is1296 -or $is1440) {
                [void]$targets.Add(@($i, $lastNumId))
            }
        }
    }
}

# --- Pass 2: apply the fix, back-to-front so earlier indices stay valid
#     while we rewrite paragraphs with InsertXML.
for ($j = $targets.Count - 1; $j -ge 0; $j--) {
    $pair = $targets[$j]
    $idx = $pair[0]
    $numId = $pair[1]

    $p = $paras.Item($idx)
    $fullXml = $p.Range.WordOpenXML

    if ($fullXml -match '(?s)<w:p\b.*?</w:p>') {
        $pxml = $Matches[0]

        # Drop the synthesized w14 paragraph/text ids - not present originally.
        $pxml = $pxml -replace ' w14:paraId="[^"]*"', ''
        $pxml = $pxml -replace ' w14:textId="[^"]*"', ''

        # WordOpenXML never reports <w:lastRenderedPageBreak/> (it's a pure
        # rendering cache hint), so InsertXML would otherwise silently drop
        # it for the one step that has it. Put it back before the run text.
        if ($pxml -notmatch 'lastRenderedPageBreak' -and $p.Range.Text -like "*Click object need sort*") {
            $pxml = $pxml -replace '(<w:r\b[^>]*><w:rPr>.*?</w:rPr>)(<w:t)', '$1<w:lastRenderedPageBreak/>$2'
        }

        $numPrXml = '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr>'

        if ($pxml -match '<w:ind w:left="1296"[^/]*/>') {
            # Step already sits at the list indent - the explicit w:ind is
            # replaced outright by the numbering properties.
            $pxml = $pxml -replace '<w:ind w:left="1296"[^/]*/>', $numPrXml
        } elseif ($pxml -match '<w:ind w:left="1440"[^/]*/>') {
            # Step used a slightly deeper manual indent - keep an explicit
            # indent but align it to the list's 1296 and add numPr before it.
            $pxml = $pxml -replace '<w:ind w:left="1440"([^/]*)/>', ($numPrXml + '<w:ind w:left="1296"$1/>')
        }

        $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body>' + $pxml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $p.Range.InsertXML($pkg)
    }
}

Write-Output ("Fixed " + $targets.Count + " step paragraph(s).")
